# Insert a new weekly record at row 98 (Terminal La Palmera de La Serena - Cebollín),
# shifting the existing rows 98..127 down to 99..128.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 98 (pushes rows 98-127 -> 99-128)
$ws.Rows.Item(98).Insert()

# Populate the new row 98 with the new weekly record.
$ws.Cells.Item(98, 1).Value2 = 8
$ws.Cells.Item(98, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(98, 3).Value2 = "Coquimbo"
$ws.Cells.Item(98, 4).Value2 = 44508
$ws.Cells.Item(98, 5).Value2 = 4
$ws.Cells.Item(98, 6).Value2 = 100112037
$ws.Cells.Item(98, 7).Value2 = "Cebollín"
$ws.Cells.Item(98, 8).Value2 = "Sin especificar"
$ws.Cells.Item(98, 9).Value2 = "Primera"
$ws.Cells.Item(98, 10).Value2 = 2800
$ws.Cells.Item(98, 11).Value2 = 900
$ws.Cells.Item(98, 12).Value2 = 1000
$ws.Cells.Item(98, 13).Value2 = 950
$ws.Cells.Item(98, 14).Value2 = "$/paquete 6 unidades"
$ws.Cells.Item(98, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(98, 16).Value2 = 158
$ws.Cells.Item(98, 17).Value2 = 6
$ws.Cells.Item(98, 18).Value2 = "Hortaliza"
